$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (Price) and volume-change (Volume(1h)) columns
# Price values are stored as text (e.g. "1.00", "25.924.01") so a leading
# apostrophe is used to force text entry and avoid Excel's automatic
# number/date conversion, which would otherwise strip formatting such as
# trailing zeros or reformat thousand-separated numbers.

$ws.Range("D2").Value = "'25.924.01"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'1.636.81"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'214.38"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -0.28%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D9").Value = "'0.0635"
$ws.Range("E9").Value = "  -0.51%  "
$ws.Range("D10").Value = "'19.48"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "'0.0795"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "'4.24"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("D13").Value = "'1.591.06"
$ws.Range("E13").Value = "  -3.03%  "
$ws.Range("E14").Value = "  -0.57%  "
$ws.Range("D15").Value = "'63.28"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "'25.965.02"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D19").Value = "'194.05"
$ws.Range("E19").Value = "  -0.04%  "
$ws.Range("E20").Value = "  -1.07%  "
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("D22").Value = "'6.17"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").Value = "'143.43"
$ws.Range("E24").Value = "  -0.48%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "  -0.32%  "
$ws.Range("E26").Value = "  -0.85%  "
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'15.51"
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -1.55%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  +0.24%  "
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +0.54%  "
$ws.Range("D35").Value = "'0.899"
$ws.Range("E35").Value = "  -0.65%  "
$ws.Range("D36").Value = "'1.125.30"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("E37").Value = "  -1.51%  "
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").Value = "'98.50"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").Value = "'0.791"
$ws.Range("E42").Value = "  -1.25%  "
$ws.Range("E43").Value = "  -0.37%  "
$ws.Range("D44").Value = "'56.26"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").Value = "'1.49"
$ws.Range("E45").Value = "  +2.23%  "
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").Value = "'7.75"
$ws.Range("E47").Value = "  +1.30%  "
$ws.Range("E48").Value = "  -0.46%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("D50").Value = "'0.0944"
$ws.Range("E50").Value = "  -2.09%  "
$ws.Range("D51").Value = "'5.49"
$ws.Range("E51").Value = "  -0.73%  "
